$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '28.539.53'
$ws.Cells.Item(2, 5).Value = '  +1.00%  '
$ws.Cells.Item(3, 4).Value = '1.570.98'
$ws.Cells.Item(3, 5).Value = '  -1.22%  '
$ws.Cells.Item(4, 5).Value = '  +0.05%  '
$origStyle = $ws.Cells.Item(5, 4).Style
$ws.Cells.Item(5, 4).Value = "'213.03"
$ws.Cells.Item(5, 4).Style = $origStyle
$ws.Cells.Item(5, 5).Value = '  -0.35%  '
$ws.Cells.Item(6, 5).Value = '  -0.40%  '
$ws.Cells.Item(7, 5).Value = '  +0.04%  '
$ws.Cells.Item(8, 5).Value = '  +4.13%  '
$origStyle = $ws.Cells.Item(9, 4).Style
$ws.Cells.Item(9, 4).Value = "'24.06"
$ws.Cells.Item(9, 4).Style = $origStyle
$ws.Cells.Item(9, 5).Value = '  -0.36%  '
$origStyle = $ws.Cells.Item(10, 4).Style
$ws.Cells.Item(10, 4).Value = "'0.247"
$ws.Cells.Item(10, 4).Style = $origStyle
$ws.Cells.Item(11, 5).Value = '  -1.56%  '
$ws.Cells.Item(12, 5).Value = '  -0.13%  '
$ws.Cells.Item(13, 4).Value = '1.795.54'
$ws.Cells.Item(13, 5).Value = '  -1.24%  '
$ws.Cells.Item(14, 4).Value = '1.568.67'
$ws.Cells.Item(14, 5).Value = '  -1.39%  '
$ws.Cells.Item(15, 5).Value = '  -1.90%  '
$ws.Cells.Item(16, 4).Value = '28.548.57'
$ws.Cells.Item(16, 5).Value = '  +0.86%  '
$origStyle = $ws.Cells.Item(17, 4).Style
$ws.Cells.Item(17, 4).Value = "'3.68"
$ws.Cells.Item(17, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(19, 4).Style
$ws.Cells.Item(19, 4).Value = "'230.15"
$ws.Cells.Item(19, 4).Style = $origStyle
$ws.Cells.Item(19, 5).Value = '  +1.00%  '
$origStyle = $ws.Cells.Item(20, 4).Style
$ws.Cells.Item(20, 4).Value = "'7.37"
$ws.Cells.Item(20, 4).Style = $origStyle
$ws.Cells.Item(20, 5).Value = '  -1.59%  '
$ws.Cells.Item(21, 5).Value = '  -2.62%  '
$ws.Cells.Item(22, 5).Value = '  +0.03%  '
$origStyle = $ws.Cells.Item(23, 4).Style
$ws.Cells.Item(23, 4).Value = "'3.88"
$ws.Cells.Item(23, 4).Style = $origStyle
$ws.Cells.Item(23, 5).Value = '  -5.51%  '
$origStyle = $ws.Cells.Item(24, 4).Style
$ws.Cells.Item(24, 4).Value = "'9.11"
$ws.Cells.Item(24, 4).Style = $origStyle
$ws.Cells.Item(24, 5).Value = '  -2.20%  '
$origStyle = $ws.Cells.Item(25, 4).Style
$ws.Cells.Item(25, 4).Value = "'2.14"
$ws.Cells.Item(25, 4).Style = $origStyle
$ws.Cells.Item(25, 5).Value = '  +9.61%  '
$origStyle = $ws.Cells.Item(26, 4).Style
$ws.Cells.Item(26, 4).Value = "'151.92"
$ws.Cells.Item(26, 4).Style = $origStyle
$ws.Cells.Item(26, 5).Value = '  +0.08%  '
$ws.Cells.Item(27, 5).Value = '  -1.15%  '
$ws.Cells.Item(28, 5).Value = '  -2.39%  '
$ws.Cells.Item(29, 5).Value = '  -3.31%  '
$ws.Cells.Item(30, 5).Value = '  +0.04%  '
$origStyle = $ws.Cells.Item(31, 4).Style
$ws.Cells.Item(31, 4).Value = "'0.0484"
$ws.Cells.Item(31, 4).Style = $origStyle
$ws.Cells.Item(31, 5).Value = '  +2.29%  '
$ws.Cells.Item(32, 5).Value = '  -2.49%  '
$ws.Cells.Item(33, 5).Value = '  -0.95%  '
$ws.Cells.Item(34, 5).Value = '  -1.73%  '
$ws.Cells.Item(35, 4).Value = '1.394.77'
$ws.Cells.Item(35, 5).Value = '  -0.45%  '
$ws.Cells.Item(36, 5).Value = '  +1.54%  '
$ws.Cells.Item(37, 5).Value = '  -3.35%  '
$ws.Cells.Item(38, 5).Value = '  +0.98%  '
$origStyle = $ws.Cells.Item(39, 4).Style
$ws.Cells.Item(39, 4).Value = "'2.62"
$ws.Cells.Item(39, 4).Style = $origStyle
$ws.Cells.Item(39, 5).Value = '  +2.96%  '
$ws.Cells.Item(40, 5).Value = '  -0.67%  '
$ws.Cells.Item(41, 5).Value = '  -2.88%  '
$ws.Cells.Item(42, 5).Value = '  +0.03%  '
$ws.Cells.Item(43, 5).Value = '  +0.80%  '
$ws.Cells.Item(44, 5).Value = '  -3.05%  '
$ws.Cells.Item(45, 5).Value = '  +2.61%  '
$ws.Cells.Item(46, 5).Value = '  -2.55%  '
$origStyle = $ws.Cells.Item(47, 4).Style
$ws.Cells.Item(47, 4).Value = "'0.968"
$ws.Cells.Item(47, 4).Style = $origStyle
$ws.Cells.Item(47, 5).Value = '  -1.85%  '
$origStyle = $ws.Cells.Item(48, 4).Style
$ws.Cells.Item(48, 4).Value = "'62.95"
$ws.Cells.Item(48, 4).Style = $origStyle
$ws.Cells.Item(48, 5).Value = '  -2.09%  '
$ws.Cells.Item(49, 4).Value = '1.708.07'
$ws.Cells.Item(49, 5).Value = '  -1.41%  '
$origStyle = $ws.Cells.Item(50, 4).Style
$ws.Cells.Item(50, 4).Value = "'86.35"
$ws.Cells.Item(50, 4).Style = $origStyle
$ws.Cells.Item(50, 5).Value = '  -1.45%  '
$ws.Cells.Item(51, 5).Value = '  -0.99%  '
